# Atualização de bases das ligas, do dia: 24-02-2024 às 12:40
#
# The rows below had their match-data columns (B, F:AC) shuffled between
# rows (the date/div columns A/C/D/E stay with their original row).
# Mapping is "destination row" -> "source row" (source row's old B..AC
# values are copied into destination row):
#   227 <- 229
#   229 <- 227
#   230 <- 233
#   231 <- 230
#   232 <- 231
#   233 <- 232
#   238 <- 240
#   239 <- 241
#   240 <- 239
#   241 <- 238

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together (everything except id/Div/DivOriginalName/Date)
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# Rows that participate in the shuffle, and where each destination row's
# data comes from (reading the *original* values before any writes).
$destToSrc = @{
    227 = 229
    229 = 227
    230 = 233
    231 = 230
    232 = 231
    233 = 232
    238 = 240
    239 = 241
    240 = 239
    241 = 238
}

# 1) Snapshot the original values of every row involved, before any
#    writes happen (several rows both supply and receive data).
$snapshot = @{}
foreach ($r in @(227,229,230,231,232,233,238,239,240,241)) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write the snapshotted source row into each destination row.
foreach ($destRow in $destToSrc.Keys) {
    $srcRow = $destToSrc[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $rowData[$c]
    }
}
